$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Range("G2").Value = 102.33
$ws.Range("G3").Value = 71.17
$ws.Range("G4").Value = 65.43000000000001
$ws.Range("G5").Value = 45.84
$ws.Range("G6").Value = 44.57
$ws.Range("G7").Value = 57.44
$ws.Range("G8").Value = 81.7
$ws.Range("G9").Value = 111.1
$ws.Range("G10").Value = 109.76
$ws.Range("G11").Value = 81.48999999999999
$ws.Range("G12").Value = 26.73
$ws.Range("G13").Value = 8.050000000000001
$ws.Range("G14").Value = 3.85
$ws.Range("G15").Value = 0.07000000000000001
$ws.Range("G16").Value = 0.02
$ws.Range("G17").Value = 1.39
$ws.Range("G18").Value = 31.96
$ws.Range("G19").Value = 88
$ws.Range("G20").Value = 111.79
$ws.Range("G21").Value = 128.94
$ws.Range("G22").Value = 140
$ws.Range("G23").Value = 142.18
$ws.Range("G24").Value = 137.05
$ws.Range("G25").Value = 121.38
